# RF007 - Gerenciar Avaliacoes : version 1.1 -> 1.2
# TC6 ("Lider de Pessoas nao e o lider" delete-error test) gains two more
# steps and is extended into a full "Editar com erro de validacao" test
# (reusing the step sequence that used to live in TC7), while TC7 becomes
# a simple happy-path delete test (reusing the pattern used by TC3/TC4/TC5).
#
# Net effect on the sheet: 3 new rows are inserted right after the old
# TC6 block (so its 4-step test grows to a 6-step test) and 3 rows are
# removed from the end of the old TC7 block (so its 6-step test shrinks
# back down to a 4-step test). Total row count / dimension is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 3 blank rows before row 72 (the empty separator
#    row that follows TC6's old 4 data rows, rows 68-71).
# ---------------------------------------------------------------------
$ws.Rows.Item(72).Insert()
$ws.Rows.Item(72).Insert()
$ws.Rows.Item(72).Insert()

# ---------------------------------------------------------------------
# 2. Update TC6's existing step 3 & step 4 (rows 70 & 71) - they change
#    from the "Excluir" flow to the "Editar" flow, matching the steps
#    that used to be TC7's steps 3 & 4.
# ---------------------------------------------------------------------
$ws.Range("B70").Value = "Lider de Pessoas clica na opcao 'Editar' para modificar a Avaliacao selecionada"
$ws.Range("D70").Value = "SYSTEM apresenta o formulario para e alteracao de Avaliacao"

$ws.Range("B71").Value = "Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' somente leitura"
$ws.Range("D71:F71").Clear()

# ---------------------------------------------------------------------
# 3. Populate the 2 newly inserted rows (73 & 74) as TC6's new step 5 &
#    step 6, matching the pattern/format of the existing numbered step
#    rows (copy format only from row 68, then set text).
# ---------------------------------------------------------------------
$ws.Range("A68:F68").Copy()
$ws.Range("A73:F73").PasteSpecial(-4122)
$ws.Range("A74:F74").PasteSpecial(-4122)

$ws.Range("A73").Value = 5
$ws.Range("B73").Value = "Lider de Pessoas seleciona 'Avaliadores' da Avaliacao"
$ws.Range("D73").Value = "SYSTEM apresenta a lista de 'Avaliadores' preenchida corretamente"

$ws.Range("A74").Value = 6
$ws.Range("B74").Value = "Lider de Pessoas clica na opcao 'Salvar'"
$ws.Range("D74").Value = "SYSTEM exibe uma mensagem de erro ao tentar editar a Avaliacao, informando o campo ou a validacao que falhou"

# ---------------------------------------------------------------------
# 4. TC7's header/sub/sub2/colheader rows (now shifted to 77-80) stay
#    exactly as they were - only its 4 data rows need updating so the
#    test becomes the "Excluir" happy-path flow.
#    Steps 1 & 2 (rows 81 & 82) are already correct (unchanged content),
#    steps 3 & 4 (rows 83 & 84, previously the "Editar" steps 3 & 4)
#    become the "Excluir" steps 3 & 4.
# ---------------------------------------------------------------------
$ws.Range("B83").Value = "Lider de Pessoas clica na opcao 'Excluir' para excluir a Avaliacao selecionada"
$ws.Range("D83").Value = "SYSTEM solicita confirmacao de exclusao da Avaliacao"

$ws.Range("B84").Value = "Lider de Pessoas confirma a exclusao do Avaliacao"

# Row 84 used to only have A/B/C filled in (it was the last of a 4-step
# block) - now it needs a D (and matching E/F) cell added, formatted the
# same way as the other "result" cells in that column.
$ws.Range("D41:F41").Copy()
$ws.Range("D84:F84").PasteSpecial(-4122)
$ws.Range("D84").Value = "SYSTEM exibe a listagem dos Avaliacoes sem o Avaliacao excluido"

# ---------------------------------------------------------------------
# 5. Remove the 3 now-superfluous rows that used to hold TC7's old
#    steps 5 & 6 plus the blank separator row between step 4 and step 5
#    (old rows 85, 86, 87 after the earlier +3 shift).
# ---------------------------------------------------------------------
$ws.Rows.Item(85).Delete()
$ws.Rows.Item(85).Delete()
$ws.Rows.Item(85).Delete()
